$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2000-2009 rows (old rows 2-11); remaining rows shift up so that
# what was 2010 (row 12) becomes row 2, etc.
$ws.Rows("2:11").Delete()

# Copy the formatting (style) of the now-shifted "2020年" row (row 12, column A)
# down into the two new rows so the new year labels get the same style (s="1").
$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)

# Add the new 2021 and 2022 rows; only column C (乡村就业人员) has data,
# columns B and D are left blank for these years.
$ws.Range("A13").Value = "2021年"
$ws.Range("C13").Value = 27879

$ws.Range("A14").Value = "2022年"
$ws.Range("C14").Value = 27420
